$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B and C are plain text (coin name / link) - safe to set directly.
# Column D values look numeric (e.g. "1.001") and column E values are
# percentage strings with spaces/signs. To guarantee these are written back
# as literal text (matching the original inlineStr cells) instead of being
# auto-converted to numbers by Excel, we force the NumberFormat of any cell
# we touch in column D to Text ("@") before assigning its value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.370.17'
$ws.Range("E2").Value = '  +1.29%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.664.00'
$ws.Range("E3").Value = '  +1.40%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.58'
$ws.Range("E5").Value = '  +1.62%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  -0.11%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3948'
$ws.Range("E7").Value = '  +0.92%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3904'
$ws.Range("E8").Value = '  +1.31%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '52.36'
$ws.Range("E9").Value = '  +6.84%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.394'
$ws.Range("E10").Value = '  +3.07%  '

$ws.Range("E11").Value = '  -0.19%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08562'
$ws.Range("E12").Value = '  +1.11%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '24.31'
$ws.Range("E13").Value = '  +1.49%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.282'
$ws.Range("E14").Value = '  +2.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.912'
$ws.Range("E15").Value = '  +5.82%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001338'
$ws.Range("E16").Value = '  +4.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.662.73'
$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '94.99'
$ws.Range("E18").Value = '  +0.78%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07000'
$ws.Range("E19").Value = '  +0.77%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.50'
$ws.Range("E20").Value = '  -1.28%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.985'
$ws.Range("E21").Value = '  +0.89%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.001'
$ws.Range("E22").Value = '  -0.18%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.72'
$ws.Range("E23").Value = '  +0.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.364.27'

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.418'
$ws.Range("E25").Value = '  +3.04%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.041'
$ws.Range("E26").Value = '  +13.32%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.51'
$ws.Range("E27").Value = '  +0.29%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.02'
$ws.Range("E28").Value = '  -0.84%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '142.34'
$ws.Range("E29").Value = '  +0.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.430'
$ws.Range("E30").Value = '  +2.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.033'
$ws.Range("E31").Value = '  -9.19%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.542'
$ws.Range("E32").Value = '  +2.95%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.844.50'
$ws.Range("E33").Value = '  +3.81%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.051'
$ws.Range("E34").Value = '  +8.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.08235'
$ws.Range("E35").Value = '  +2.21%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.03017'
$ws.Range("E36").Value = '  +3.05%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.895'
$ws.Range("E37").Value = '  -3.44%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '11.11'
$ws.Range("E38").Value = '  +11.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2750'
$ws.Range("E39").Value = '  +1.80%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09223'
$ws.Range("E40").Value = '  -0.24%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7686'
$ws.Range("E41").Value = '  +0.92%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.76'
$ws.Range("E42").Value = '  +5.15%  '

$ws.Range("E43").Value = '  -1.35%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.57'
$ws.Range("E44").Value = '  +3.88%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.7085'
$ws.Range("E45").Value = '  +2.99%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.527'
$ws.Range("E46").Value = '  +1.91%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.129'
$ws.Range("E47").Value = '  +1.06%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  -0.11%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08409'
$ws.Range("E49").Value = '  +0.39%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.70'
$ws.Range("E50").Value = '  +2.09%  '

$ws.Range("B51").Value = 'Tezos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fsIbGOEJWbzxG+tezos-xtz'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.428'
$ws.Range("E51").Value = '  +11.98%  '
